$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# taskKey (column C) for the "e.init" activity rows was renamed to "e.adjust"
$ws.Range("C21:C26").Value = "e.adjust"

# Match the author's final cursor/selection position recorded in the workbook
$ws.Range("E23").Select()
